$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sheet view stays left-to-right (matches the explicit rightToLeft="false"
# now recorded on the sheetView).
$ws.DisplayRightToLeft = $false

# --- Update existing row 2 values ---
$ws.Range("B2").Value = "Name Kumar1"
$ws.Range("G2").Value = "5/9/2019 9:40:00 PM"

# --- Add new row 3 values ---
$ws.Range("A3").Value = "49cd1269-d104-4a11-9985-561a807f6c64"
$ws.Range("B3").Value = "Product2"
$ws.Range("C3").Value = "SKU2"
$ws.Range("D3").Value = "Code2"
$ws.Range("E3").Value = "Desc2"
$ws.Range("F3").Value = "Active"
$ws.Range("G3").Value = "5/9/2019 9:40:28 PM"

# --- Re-apply row 2's existing cell style (quote-prefixed "Normal") to every
#     cell in rows 2 and 3, since assigning .Value resets a cell's style. We
#     use A2 (untouched, still carrying the original style) as the format
#     donor. xlPasteFormats (-4122) copies formatting only, not values. ---
$ws.Range("A2").Copy()
$ws.Range("A2:G3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
